$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.789.25"
$ws.Range("E2").Value = '  +6.90%  '

$ws.Range("D3").Value = "'1.761.12"
$ws.Range("E3").Value = '  +5.43%  '

$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").Value = "'316.81"
$ws.Range("E5").Value = '  +3.02%  '

$ws.Range("E6").Value = '  +0.45%  '

$ws.Range("D7").Value = "'0.3825"
$ws.Range("E7").Value = '  +3.18%  '

$ws.Range("D8").Value = "'0.3608"
$ws.Range("E8").Value = '  +5.06%  '

$ws.Range("D9").Value = "'50.25"
$ws.Range("E9").Value = '  +4.65%  '

$ws.Range("D10").Value = "'1.227"
$ws.Range("E10").Value = '  +4.49%  '

$ws.Range("D11").Value = "'0.07707"
$ws.Range("E11").Value = '  +6.35%  '

$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("E13").Value = '  +5.78%  '

$ws.Range("D14").Value = "'6.478"
$ws.Range("E14").Value = '  +7.54%  '

$ws.Range("D15").Value = "'7.100"
$ws.Range("E15").Value = '  +5.24%  '

$ws.Range("D16").Value = "'1.763.25"
$ws.Range("E16").Value = '  +5.78%  '

$ws.Range("D17").Value = "'0.00001157"

$ws.Range("D18").Value = "'0.06787"
$ws.Range("E18").Value = '  +1.25%  '

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").Value = "'86.91"
$ws.Range("E20").Value = '  +6.45%  '

$ws.Range("D21").Value = "'17.70"
$ws.Range("E21").Value = '  +7.75%  '

$ws.Range("D22").Value = "'6.523"
$ws.Range("E22").Value = '  +5.89%  '

$ws.Range("D23").Value = "'13.03"
$ws.Range("E23").Value = '  +8.66%  '

$ws.Range("D24").Value = "'25.759.84"
$ws.Range("E24").Value = '  +7.10%  '

$ws.Range("D25").Value = "'2.442"
$ws.Range("E25").Value = '  +1.73%  '

$ws.Range("D26").Value = "'2.910"
$ws.Range("E26").Value = '  +9.06%  '

$ws.Range("D27").Value = "'20.84"
$ws.Range("E27").Value = '  +6.77%  '

$ws.Range("D28").Value = "'155.94"
$ws.Range("E28").Value = '  +2.62%  '

$ws.Range("D29").Value = "'1.959.88"
$ws.Range("E29").Value = '  +5.83%  '

$ws.Range("D30").Value = "'133.91"
$ws.Range("E30").Value = '  +5.42%  '

$ws.Range("D31").Value = "'1.209"
$ws.Range("E31").Value = '  +23.12%  '

$ws.Range("D32").Value = "'7.246"
$ws.Range("E32").Value = '  +14.42%  '

$ws.Range("D33").Value = "'4.213"
$ws.Range("E33").Value = '  +3.56%  '

$ws.Range("D34").Value = "'14.37"
$ws.Range("E34").Value = '  +16.71%  '

$ws.Range("D35").Value = "'1.813"
$ws.Range("E35").Value = '  +4.91%  '

$ws.Range("D36").Value = "'0.08766"
$ws.Range("E36").Value = '  +4.64%  '

$ws.Range("D37").Value = "'5.740"
$ws.Range("E37").Value = '  +7.71%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = "'9.416"
$ws.Range("E38").Value = '  +5.13%  '

$ws.Range("D39").Value = "'0.06759"
$ws.Range("E39").Value = '  +6.28%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "'0.02501"
$ws.Range("E40").Value = '  +7.38%  '

$ws.Range("D41").Value = "'0.2267"
$ws.Range("E41").Value = '  +8.93%  '

$ws.Range("D42").Value = "'1.301"
$ws.Range("E42").Value = '  +0.93%  '

$ws.Range("D43").Value = "'0.6604"
$ws.Range("E43").Value = '  +8.24%  '

$ws.Range("D44").Value = "'14.37"
$ws.Range("E44").Value = '  +9.89%  '

$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = '  +0.53%  '

$ws.Range("D46").Value = "'0.6363"
$ws.Range("E46").Value = '  +7.09%  '

$ws.Range("D47").Value = "'3.900"
$ws.Range("E47").Value = '  +2.39%  '

$ws.Range("D48").Value = "'2.179"
$ws.Range("E48").Value = '  +8.68%  '

$ws.Range("D49").Value = "'132.15"
$ws.Range("E49").Value = '  +3.81%  '

$ws.Range("D50").Value = "'0.07503"
$ws.Range("E50").Value = '  +5.49%  '

$ws.Range("D51").Value = "'81.12"
$ws.Range("E51").Value = '  +6.97%  '
